# Add two new report-export columns ("code" and "nextaction") right before
# the existing "nextkin" / "kinphone" columns (which shift from J/K to L/M),
# extend the conditional-formatting range and duplicate-value highlighting
# over the new columns, and move the active selection to the newly added
# area - matching the "added new files of report expoerts" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at J (old J "nextkin" / K "kinphone" shift to L/M).
$ws.Range("J1").EntireColumn.Insert()
$ws.Range("J1").EntireColumn.Insert()

# New headers for the inserted columns.
$ws.Range("J1").Value = "code"
$ws.Range("K1").Value = "nextaction"

# Grow the "duplicate values" conditional formatting that used to cover
# H2:J14 so it also covers the two new columns (now H2:L14).
$dupRule = $ws.Range("H2:J14").FormatConditions.Item(1)
$dupRule.ModifyAppliesToRange($ws.Range("H2:L14"))

# Leave the selection on the newly added area, like in the saved workbook.
$ws.Range("L5").Select()
